$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Headers: BTec_Logo-Orange pictures currently named "image2.jpg" -> rename to "image1.jpg"
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hdr = $sec.Headers($i)
    if ($hdr.Range.InlineShapes.Count -gt 0) {
        $shp = $hdr.Range.InlineShapes(1)
        if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
            $shp.Name = "image1.jpg"
            Write-Host "Renamed header logo $i to image1.jpg"
        }
    }
}

# Footers: PearsonLogo pictures currently named "image1.png" -> rename to "image2.png"
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers($i)
    if ($ftr.Range.InlineShapes.Count -gt 0) {
        $shp = $ftr.Range.InlineShapes(1)
        if ($shp.AlternativeText -like "*PearsonLogo.png") {
            $shp.Name = "image2.png"
            Write-Host "Renamed footer logo $i to image2.png"
        }
    }
}
